$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantity (F) and value (G) for affected stock rows
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("F37").Value = 18
$ws.Range("G37").Value = 829.8
$ws.Range("F47").Value = 24
$ws.Range("G47").Value = 873.84
$ws.Range("F66").Value = 52
$ws.Range("G66").Value = 2319.2
$ws.Range("F76").Value = 60
$ws.Range("G76").Value = 2007.6
$ws.Range("F138").Value = 6
$ws.Range("G138").Value = 185.16
$ws.Range("F140").Value = 43
$ws.Range("G140").Value = 1824.49
$ws.Range("F145").Value = 2
$ws.Range("G145").Value = 64.28
$ws.Range("F146").Value = 71
$ws.Range("G146").Value = 1368.88
$ws.Range("F162").Value = 205
$ws.Range("G162").Value = 4034.4
$ws.Range("F171").Value = 7
$ws.Range("G171").Value = 457.1
$ws.Range("F174").Value = 11
$ws.Range("G174").Value = 1021.79
$ws.Range("F230").Value = 12
$ws.Range("G230").Value = 1521.48
$ws.Range("F233").Value = 81
$ws.Range("G233").Value = 1513.89
$ws.Range("F251").Value = 7
$ws.Range("G251").Value = 1861.72
$ws.Range("F290").Value = 31
$ws.Range("G290").Value = 4319.85
$ws.Range("F340").Value = 13
$ws.Range("G340").Value = 958.23
$ws.Range("F351").Value = 56
$ws.Range("G351").Value = 1294.72
$ws.Range("F356").Value = 1
$ws.Range("G356").Value = 135.55
$ws.Range("F359").Value = 4
$ws.Range("G359").Value = 685.32
$ws.Range("F361").Value = 155
$ws.Range("G361").Value = 6382.9
$ws.Range("F362").Value = 32
$ws.Range("G362").Value = 2964.16
$ws.Range("F363").Value = 310
$ws.Range("G363").Value = 12322.5
$ws.Range("F366").Value = 2
$ws.Range("G366").Value = 306.82
$ws.Range("F382").Value = 11
$ws.Range("G382").Value = 857.12
$ws.Range("F388").Value = 1
$ws.Range("G388").Value = 42.01
$ws.Range("F392").Value = 220
$ws.Range("G392").Value = 4485.8
$ws.Range("F406").Value = 43
$ws.Range("G406").Value = 6461.61
$ws.Range("F433").Value = 618
$ws.Range("G433").Value = 59698.8
$ws.Range("F436").Value = 112
$ws.Range("G436").Value = 3011.68
$ws.Range("F455").Value = 7
$ws.Range("G455").Value = 201.46
$ws.Range("F461").Value = 450
$ws.Range("G461").Value = 5805
$ws.Range("F463").Value = 81
$ws.Range("G463").Value = 4094.55
$ws.Range("F464").Value = 641
$ws.Range("G464").Value = 8268.9
$ws.Range("F481").Value = 568
$ws.Range("G481").Value = 14938.4
$ws.Range("F485").Value = 474
$ws.Range("G485").Value = 9352.02
$ws.Range("F488").Value = 275
$ws.Range("G488").Value = 5351.5
$ws.Range("F491").Value = 1425
$ws.Range("G491").Value = 9248.25
$ws.Range("F494").Value = 477
$ws.Range("G494").Value = 7837.11
$ws.Range("F499").Value = 4
$ws.Range("G499").Value = 125.2
$ws.Range("F501").Value = 16
$ws.Range("G501").Value = 801.28
$ws.Range("F502").Value = 143
$ws.Range("G502").Value = 4149.86
$ws.Range("F503").Value = 56
$ws.Range("G503").Value = 3452.4
$ws.Range("F505").Value = 17
$ws.Range("G505").Value = 548.59
$ws.Range("F514").Value = 9
$ws.Range("G514").Value = 884.07
$ws.Range("F520").Value = 16
$ws.Range("G520").Value = 5513.12
$ws.Range("F521").Value = 556
$ws.Range("G521").Value = 5960.32
$ws.Range("F524").Value = 104
$ws.Range("G524").Value = 6314.88
$ws.Range("F605").Value = 52
$ws.Range("G605").Value = 1864.72
$ws.Range("F650").Value = 54
$ws.Range("G650").Value = 4483.08
$ws.Range("F721").Value = 51
$ws.Range("G721").Value = 7156.83
$ws.Range("F750").Value = 159
$ws.Range("G750").Value = 17734.86
$ws.Range("F758").Value = 86
$ws.Range("G758").Value = 4620.78
$ws.Range("F759").Value = 172
$ws.Range("G759").Value = 25994.36
$ws.Range("F780").Value = 57
$ws.Range("G780").Value = 2556.45
$ws.Range("F786").Value = 330
$ws.Range("G786").Value = 25934.7
$ws.Range("F787").Value = 21
$ws.Range("G787").Value = 7934.01
$ws.Range("F790").Value = 271
$ws.Range("G790").Value = 9980.93

# Update Sub Total / Grand Total (B) rows
$ws.Range("B15").Value = 27128.76
$ws.Range("B41").Value = 21313.14
$ws.Range("B77").Value = 128109.86
$ws.Range("B147").Value = 57459.72
$ws.Range("B168").Value = 29040.66
$ws.Range("B183").Value = 28832.39
$ws.Range("B237").Value = 13433.96
$ws.Range("B264").Value = 88725.47
$ws.Range("B336").Value = 252725.88
$ws.Range("B364").Value = 53816.4
$ws.Range("B370").Value = 3438.74
$ws.Range("B397").Value = 24113.28
$ws.Range("B407").Value = 14326.83
$ws.Range("B439").Value = 130750.5
$ws.Range("B459").Value = 11548.82
$ws.Range("B470").Value = 39055.8
$ws.Range("B496").Value = 164778.37
$ws.Range("B515").Value = 37606.26
$ws.Range("B539").Value = 115696.21
$ws.Range("B625").Value = 45858.46
$ws.Range("B651").Value = 241737.81
$ws.Range("B722").Value = 27393.64
$ws.Range("B767").Value = 566324.62
$ws.Range("B796").Value = 163074.65
$ws.Range("B855").Value = 5294574.19
$ws.Range("B856").Value = 5294574.19
